$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin rows scraped by the refresh job: cell address -> new text value.
# (Matches the 'Updated cryptos list ... with GitHub Actions' price/volume refresh.)
$updates = [ordered]@{
    'D2' = '61.977.36'
    'E2' = '  +0.81%  '
    'D3' = '3.416.31'
    'E3' = '  +1.18%  '
    'E4' = '  -0.04%  '
    'D5' = '410.11'
    'E5' = '  +0.82%  '
    'D6' = '128.76'
    'E6' = '  -4.69%  '
    'D7' = '0.624'
    'E7' = '  +5.10%  '
    'E8' = '  -0.11%  '
    'D9' = '0.752'
    'E9' = '  +12.16%  '
    'D10' = '0.140'
    'E10' = '  +16.03%  '
    'D11' = '42.93'
    'E11' = '  +0.65%  '
    'E12' = '  -0.46%  '
    'D13' = '21.27'
    'E13' = '  +7.89%  '
    'D14' = '8.86'
    'E14' = '  +5.21%  '
    'D15' = '0.0000203'
    'E15' = '  +59.08%  '
    'D16' = '3.449.67'
    'E16' = '  +1.77%  '
    'D17' = '12.61'
    'E17' = '  +14.65%  '
    'E18' = '  +4.12%  '
    'D19' = '61.945.10'
    'E19' = '  +0.84%  '
    'D20' = '406.35'
    'E20' = '  +29.19%  '
    'D21' = '91.01'
    'E21' = '  +6.98%  '
    'E22' = '  -0.64%  '
    'D23' = '13.41'
    'E23' = '  +4.53%  '
    'D24' = '3.25'
    'E24' = '  +3.16%  '
    'D25' = '33.06'
    'E25' = '  +11.79%  '
    'D26' = '4.80'
    'E26' = '  +0.17%  '
    'D27' = '8.53'
    'E27' = '  +1.71%  '
    'D28' = '7.63'
    'E28' = '  -0.33%  '
    'D29' = '2.71'
    'E29' = '  +4.72%  '
    'E30' = '  +0.91%  '
    'E31' = '  +0.54%  '
    'D32' = '43.96'
    'E32' = '  +8.15%  '
    'D33' = '11.81'
    'E33' = '  +3.97%  '
    'D34' = '0.999'
    'E34' = '  +0.04%  '
    'D35' = '0.0500'
    'E35' = '  +3.74%  '
    'D36' = '52.87'
    'E36' = '  +1.84%  '
    'E37' = '  +0.00%  '
    'E38' = '  -0.78%  '
    'B39' = 'Stellar'
    'C39' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D39' = '0.132'
    'E39' = '  +6.56%  '
    'B40' = 'Stacks'
    'C40' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D40' = '2.90'
    'E40' = '  -1.00%  '
    'D41' = '0.315'
    'E41' = '  +6.37%  '
    'D42' = '140.71'
    'E42' = '  +1.56%  '
    'E43' = '  -0.17%  '
    'D44' = '4.01'
    'E44' = '  -0.77%  '
    'B45' = 'Celestia'
    'C45' = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
    'D45' = '16.81'
    'E45' = '  +0.29%  '
    'B46' = 'WEMIXToken'
    'C46' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D46' = '2.36'
    'E46' = '  +6.22%  '
    'D47' = '21.89'
    'E47' = '  +2.69%  '
    'D48' = '2.108.56'
    'E48' = '  -0.71%  '
    'B49' = 'ApeXProtocol'
    'C49' = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    'D49' = '2.29'
    'E49' = '  -0.25%  '
    'D50' = '1.93'
    'E50' = '  +0.16%  '
    'B51' = 'BEAM'
    'C51' = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
    'D51' = '0.0373'
    'E51' = '  +7.85%  '
}

# Column D holds prices as plain text (e.g. '61.977.36', '0.0373') so that thousand-
# separator dots and trailing zeros survive. Several of the new values parse as valid
# numbers (e.g. '0.624', '42.93', '0.140'), so without forcing Text format first Excel
# would silently convert them to numbers and drop meaningful trailing zeros.
$textForceCells = @(
    'D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D13', 'D14', 'D15', 'D17', 'D20', 'D21',
    'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D32', 'D33', 'D34', 'D35', 'D36',
    'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50', 'D51'
)

foreach ($cell in $textForceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

# Restore the default style on the text-forced cells so no stray per-cell NumberFormat
# is left behind once the text value has been committed.
foreach ($cell in $textForceCells) {
    $ws.Range($cell).Style = "Normal"
}

